$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings
# (e.g. "1.001", "5.281") are preserved as text, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '30.270.94'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '1.862.48'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '236.29'
$ws.Range('E5').Value = '  +0.38%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').Value = '0.4705'
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('D8').Value = '0.2915'
$ws.Range('E8').Value = '  +2.36%  '
$ws.Range('D9').Value = '0.06549'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').Value = '21.89'
$ws.Range('E10').Value = '  +3.00%  '
$ws.Range('D11').Value = '0.07931'
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('D12').Value = '97.90'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').Value = '1.862.51'
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').Value = '5.147'
$ws.Range('E14').Value = '  +0.54%  '
$ws.Range('D15').Value = '0.6806'
$ws.Range('E15').Value = '  +0.55%  '
$ws.Range('D16').Value = '264.37'
$ws.Range('E16').Value = '  -6.22%  '
$ws.Range('D17').Value = '30.247.79'
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '13.82'
$ws.Range('E18').Value = '  +8.98%  '
$ws.Range('D19').Value = '1.0000'
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('D20').Value = '0.000007452'
$ws.Range('E20').Value = '  +2.12%  '
$ws.Range('D21').Value = '2.108.04'
$ws.Range('E21').Value = '  -0.27%  '
$ws.Range('B22').Value = 'BinanceUSD'
$ws.Range('C22').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D22').Value = '1.000'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '5.281'
$ws.Range('E23').Value = '  -4.16%  '
$ws.Range('D24').Value = '6.176'
$ws.Range('E24').Value = '  -0.01%  '
$ws.Range('D25').Value = '167.50'
$ws.Range('E25').Value = '  +1.69%  '
$ws.Range('D26').Value = '9.206'
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('D27').Value = '18.89'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').Value = '1.954'
$ws.Range('E28').Value = '  +1.21%  '
$ws.Range('D29').Value = '1.394'
$ws.Range('E29').Value = '  +1.50%  '
$ws.Range('D30').Value = '0.09858'
$ws.Range('E30').Value = '  +1.41%  '
$ws.Range('E31').Value = '  -1.39%  '
$ws.Range('D32').Value = '1.470'
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('D33').Value = '4.039'
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('E34').Value = '  +0.35%  '
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('D36').Value = '0.7032'
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('D38').Value = '0.01880'
$ws.Range('E38').Value = '  +0.88%  '
$ws.Range('D39').Value = '2.619'
$ws.Range('E39').Value = '  +2.97%  '
$ws.Range('D40').Value = '6.312'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').Value = '74.01'
$ws.Range('E41').Value = '  +1.05%  '
$ws.Range('D42').Value = '1.949'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = '0.8480'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = '0.4161'
$ws.Range('E44').Value = '  -0.57%  '
$ws.Range('D45').Value = '0.9987'
$ws.Range('E45').Value = '  -0.19%  '
$ws.Range('D46').Value = '103.27'
$ws.Range('E46').Value = '  -0.62%  '
$ws.Range('D47').Value = '7.170'
$ws.Range('E47').Value = '  -0.45%  '
$ws.Range('D48').Value = '944.42'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('D49').Value = '9.183'
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('E50').Value = '  +0.24%  '
$ws.Range('D51').Value = '0.05662'
$ws.Range('E51').Value = '  +0.61%  '

# Restore default style on column D so no stray number-format style remains
$ws.Range("D2:D51").Style = "Normal"
